# "Generate Report for Handback"
#
# The handback run for both locales (zh-cn, de-de) completed: the report's
# Status column moves from "In Translation" to "Handed back: in sync with
# en-US", the "Latest Target File" / "Latest Handback File" / "Latest
# Handback DateTime" columns get populated with the handback artifacts, and
# a couple of columns are widened so the new (longer) values are readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1) Status: "In Translation" -> "Handed back: in sync with en-US"
#    (Overview mirrors the per-locale Status in its zh-cn / de-de columns)
# ---------------------------------------------------------------------
$handedBack = "Handed back: in sync with en-US"

$overview.Range("E2").Value = $handedBack
$overview.Range("F2").Value = $handedBack
$overview.Range("E3").Value = $handedBack
$overview.Range("F3").Value = $handedBack

$zhcn.Range("C2").Value = $handedBack
$zhcn.Range("C3").Value = $handedBack

$dede.Range("C2").Value = $handedBack
$dede.Range("C3").Value = $handedBack

# ---------------------------------------------------------------------
# 2) Populate "Latest Target File" (I), "Latest Handback File" (J) and
#    "Latest Handback DateTime" (K) for both rows on both locale sheets.
# ---------------------------------------------------------------------

# -- zh-cn --
$zhcn.Range("J2").Value = "588fbcad-436d-4462-8ae0-6ce9bf9bc1f6.fb2c1b3b466b7f102d42f7ebbb8a27de6011cec1.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-16 18:20:46"

$zhcn.Range("J3").Value = "ac74fad0-c4a7-4c40-ac88-bc6c0f79e9e4.7c9b6143a63de52480f3a22aefc9b953b02372e7.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-16 18:20:46"

# -- de-de --
$dede.Range("J2").Value = "588fbcad-436d-4462-8ae0-6ce9bf9bc1f6.fb2c1b3b466b7f102d42f7ebbb8a27de6011cec1.de-de.xlf"
$dede.Range("K2").Value = "2016-08-16 18:20:52"

$dede.Range("J3").Value = "ac74fad0-c4a7-4c40-ac88-bc6c0f79e9e4.7c9b6143a63de52480f3a22aefc9b953b02372e7.de-de.xlf"
$dede.Range("K3").Value = "2016-08-16 18:20:52"

# ---------------------------------------------------------------------
# 3) "Latest Target File" (I2/I3) becomes a hyperlink to the source .md,
#    same as column A -- so rebuild each sheet's hyperlinks, preserving
#    the existing A2/A3 links and inserting the new I2/I3 ones in place.
# ---------------------------------------------------------------------
function Rebuild-Hyperlinks($ws) {
    $urlA2 = $null
    $urlA3 = $null
    foreach ($hl in $ws.Hyperlinks) {
        $refAddr = $hl.Range.Address()
        if ($refAddr -eq '$A$2') { $urlA2 = $hl.Address }
        if ($refAddr -eq '$A$3') { $urlA3 = $hl.Address }
    }

    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $urlA2, "", "", "588fbcad-436d-4462-8ae0-6ce9bf9bc1f6.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlA2, "", "", "588fbcad-436d-4462-8ae0-6ce9bf9bc1f6.md")
    $ws.Range("I2").Style = "HyperLink"

    $ws.Hyperlinks.Add($ws.Range("A3"), $urlA3, "", "", "ac74fad0-c4a7-4c40-ac88-bc6c0f79e9e4.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlA3, "", "", "ac74fad0-c4a7-4c40-ac88-bc6c0f79e9e4.md")
    $ws.Range("I3").Style = "HyperLink"
}

Rebuild-Hyperlinks($zhcn)
Rebuild-Hyperlinks($dede)

# ---------------------------------------------------------------------
# 4) Widen columns to fit the newly-populated / lengthened values.
#    (ColumnWidth is in characters; these were solved empirically so the
#    saved XML column width lands on/near the target.)
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.1   # E: zh-cn
$overview.Columns.Item(6).ColumnWidth = 29.1   # F: de-de

$zhcn.Columns.Item(3).ColumnWidth = 29.1    # C: Status
$zhcn.Columns.Item(9).ColumnWidth = 39.17   # I: Latest Target File
$zhcn.Columns.Item(10).ColumnWidth = 39.17  # J: Latest Handback File

$dede.Columns.Item(3).ColumnWidth = 29.1    # C: Status
$dede.Columns.Item(9).ColumnWidth = 39.17   # I: Latest Target File
$dede.Columns.Item(10).ColumnWidth = 39.17  # J: Latest Handback File
